# Generate Report for Handback
# -----------------------------------------------------------------------------
# This script mirrors what the Open Localization handback bot does after a
# translation round-trip completes: it marks the overall status as handed
# back (now in sync with en-US), records the freshly produced target/handback
# files for each locale, stamps the handback datetime, and adds a hyperlink
# to the newly produced target file - same as the already-existing hyperlink
# on the source file name column.

$wb = $excel.ActiveWorkbook

$hyperlinkUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f76a5d86231675c720c5f939ba6b3642d229df69/e2e/b85bd77a-a246-4ec4-b1d1-799fc8ead735.md"
$targetDisplay = "b85bd77a-a246-4ec4-b1d1-799fc8ead735.md"
$handedBackStatus = "Handed back: in sync with en-US"

# --- Overview sheet -----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $handedBackStatus
$wsOverview.Range("F2").Value = $handedBackStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $handedBackStatus

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, $targetDisplay)
$wsZhCn.Range("J2").Value = "b85bd77a-a246-4ec4-b1d1-799fc8ead735.4d6f70bac370f301352a11648ffc839c2aa4d5e2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-10-18 04:52:33"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $handedBackStatus

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, $targetDisplay)
$wsDeDe.Range("J2").Value = "b85bd77a-a246-4ec4-b1d1-799fc8ead735.4d6f70bac370f301352a11648ffc839c2aa4d5e2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-10-18 04:52:56"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
